# Auto-generated Excel COM-interop script to apply the Ravana_Profits.xlsx diff.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific
# rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 80
$ws.Range("H80").Value = 712.5
$ws.Range("I80").Value = 712.5
$ws.Range("K80").Value = 2137.5
$ws.Range("M80").Value = -1139.5

# ALC row 83
$ws.Range("H83").Value = 712.5
$ws.Range("I83").Value = 712.5
$ws.Range("K83").Value = 6412.5
$ws.Range("M83").Value = -1420.5

# ALC row 94
$ws.Range("H94").Value = 3996.25
$ws.Range("I94").Value = 3996.25
$ws.Range("K94").Value = 3996.25
$ws.Range("M94").Value = -3545.25

# ALC row 96
$ws.Range("H96").Value = 857.125
$ws.Range("I96").Value = 611.6
$ws.Range("J96").Value = 1266.3334
$ws.Range("K96").Value = 1834.8
$ws.Range("L96").Value = 3799.0002
$ws.Range("M96").Value = -461.8000000000002
$ws.Range("N96").Value = -6545.0002

# ALC row 98
$ws.Range("H98").Value = 738.3333
$ws.Range("I98").Value = 335.14285
$ws.Range("K98").Value = 335.14285
$ws.Range("M98").Value = 1162.85715

# ALC row 122
$ws.Range("H122").Value = 738.3333
$ws.Range("I122").Value = 335.14285
$ws.Range("K122").Value = 1005.42855
$ws.Range("M122").Value = 1444.57145

# ALC row 132
$ws.Range("H132").Value = 1305.55
$ws.Range("I132").Value = 1173.4445
$ws.Range("J132").Value = 2494.5
$ws.Range("K132").Value = 3520.3335
$ws.Range("L132").Value = 7483.5
$ws.Range("M132").Value = -990.3335000000002
$ws.Range("N132").Value = -12543.5

# ALC row 137
$ws.Range("H137").Value = 3063.3809
$ws.Range("I137").Value = 917.625
$ws.Range("J137").Value = 4383.846
$ws.Range("K137").Value = 2752.875
$ws.Range("L137").Value = 13151.538
$ws.Range("M137").Value = -202.875
$ws.Range("N137").Value = -18251.538

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 1123.8572
$ws.Range("I2").Value = 978.6667
$ws.Range("J2").Value = 1995
$ws.Range("K2").Value = 978.6667
$ws.Range("L2").Value = 1995
$ws.Range("M2").Value = -865.6667
$ws.Range("N2").Value = -2221

# ARM row 102
$ws.Range("H102").Value = 2077
$ws.Range("I102").Value = 1729.6666
$ws.Range("K102").Value = 1729.6666
$ws.Range("M102").Value = -107.6666

# ARM row 116
$ws.Range("H116").Value = 1123.8572
$ws.Range("I116").Value = 978.6667
$ws.Range("J116").Value = 1995
$ws.Range("K116").Value = 978.6667
$ws.Range("L116").Value = 1995
$ws.Range("M116").Value = 1315.3333
$ws.Range("N116").Value = -6583

# ARM row 122
$ws.Range("H122").Value = 1569.5454
$ws.Range("I122").Value = 1604.3
$ws.Range("J122").Value = 1222
$ws.Range("K122").Value = 4812.9
$ws.Range("L122").Value = 3666
$ws.Range("M122").Value = -2362.9
$ws.Range("N122").Value = -8566

# ARM row 132
$ws.Range("H132").Value = 2289.762
$ws.Range("I132").Value = 1425.3462
$ws.Range("K132").Value = 4276.0386
$ws.Range("M132").Value = -1746.0386

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 1123.8572
$ws.Range("I3").Value = 978.6667
$ws.Range("J3").Value = 1995
$ws.Range("K3").Value = 978.6667
$ws.Range("L3").Value = 1995
$ws.Range("M3").Value = -864.6667
$ws.Range("N3").Value = -2223

# BSM row 20
$ws.Range("H20").Value = 2472
$ws.Range("I20").Value = 2472
$ws.Range("K20").Value = 2472
$ws.Range("M20").Value = -2225

# BSM row 86
$ws.Range("H86").Value = 4517.706
$ws.Range("I86").Value = 3902.5
$ws.Range("J86").Value = 5064.5557
$ws.Range("K86").Value = 3902.5
$ws.Range("L86").Value = 5064.5557
$ws.Range("M86").Value = -2779.5
$ws.Range("N86").Value = -7310.5557

# BSM row 89
$ws.Range("H89").Value = 4517.706
$ws.Range("I89").Value = 3902.5
$ws.Range("J89").Value = 5064.5557
$ws.Range("K89").Value = 19512.5
$ws.Range("L89").Value = 25322.7785
$ws.Range("M89").Value = -13896.5
$ws.Range("N89").Value = -36554.7785

# BSM row 94
$ws.Range("H94").Value = 505.375
$ws.Range("I94").Value = 468.16666
$ws.Range("J94").Value = 617
$ws.Range("K94").Value = 468.16666
$ws.Range("L94").Value = 617
$ws.Range("M94").Value = -17.16665999999998
$ws.Range("N94").Value = -1519

# BSM row 99
$ws.Range("H99").Value = 1666.3334
$ws.Range("J99").Value = 1499
$ws.Range("L99").Value = 1499
$ws.Range("N99").Value = -4495

$ws = $wb.Worksheets.Item("CRP")
# CRP row 7
$ws.Range("H7").Value = 345.7
$ws.Range("I7").Value = 270
$ws.Range("J7").Value = 396.16666
$ws.Range("K7").Value = 270
$ws.Range("L7").Value = 396.16666
$ws.Range("M7").Value = -157
$ws.Range("N7").Value = -622.16666

# CRP row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = ""

# CRP row 99
$ws.Range("H99").Value = 4904.1665
$ws.Range("I99").Value = 4846.6665
$ws.Range("J99").Value = 4961.6665
$ws.Range("K99").Value = 4846.6665
$ws.Range("L99").Value = 4961.6665
$ws.Range("M99").Value = -3348.6665
$ws.Range("N99").Value = -7957.6665

# CRP row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""

# CRP row 126
$ws.Range("H126").Value = 4904.1665
$ws.Range("I126").Value = 4846.6665
$ws.Range("J126").Value = 4961.6665
$ws.Range("K126").Value = 14539.9995
$ws.Range("L126").Value = 14884.9995
$ws.Range("M126").Value = -12069.9995
$ws.Range("N126").Value = -19824.9995

$ws = $wb.Worksheets.Item("CUL")
# CUL row 80
$ws.Range("H80").Value = 502195.5
$ws.Range("I80").Value = 4392
$ws.Range("K80").Value = 13176
$ws.Range("M80").Value = -12240

# CUL row 83
$ws.Range("H83").Value = 502195.5
$ws.Range("I83").Value = 4392
$ws.Range("K83").Value = 39528
$ws.Range("M83").Value = -34848

# CUL row 94
$ws.Range("H94").Value = 1410
$ws.Range("I94").Value = 1410
$ws.Range("K94").Value = 4230
$ws.Range("M94").Value = -3554

# CUL row 107
$ws.Range("H107").Value = 554.8889
$ws.Range("J107").Value = 554.8889
$ws.Range("L107").Value = 1664.6667
$ws.Range("N107").Value = -5504.6667

# CUL row 132
$ws.Range("H132").Value = 4122.125
$ws.Range("I132").Value = 4245
$ws.Range("J132").Value = 3999.25
$ws.Range("K132").Value = 38205
$ws.Range("L132").Value = 35993.25
$ws.Range("M132").Value = -35675
$ws.Range("N132").Value = -41053.25

# CUL row 134
$ws.Range("H134").Value = 2032
$ws.Range("I134").Value = 2032
$ws.Range("K134").Value = 6096
$ws.Range("M134").Value = -1026

$ws = $wb.Worksheets.Item("GSM")
# GSM row 102
$ws.Range("H102").Value = 2738.8235
$ws.Range("I102").Value = 2036.4
$ws.Range("J102").Value = 8007
$ws.Range("K102").Value = 2036.4
$ws.Range("L102").Value = 8007
$ws.Range("M102").Value = -414.4000000000001
$ws.Range("N102").Value = -11251

# GSM row 107
$ws.Range("H107").Value = 1048.5
$ws.Range("I107").Value = 1065.1666
$ws.Range("J107").Value = 998.5
$ws.Range("K107").Value = 1065.1666
$ws.Range("L107").Value = 998.5
$ws.Range("M107").Value = 854.8334
$ws.Range("N107").Value = -4838.5

# GSM row 126
$ws.Range("H126").Value = 945.2
$ws.Range("I126").Value = 944.6667
$ws.Range("J126").Value = 946
$ws.Range("K126").Value = 2834.0001
$ws.Range("L126").Value = 2838
$ws.Range("M126").Value = -364.0001000000002
$ws.Range("N126").Value = -7778

# GSM row 132
$ws.Range("H132").Value = 1758.7858
$ws.Range("I132").Value = 1207.65
$ws.Range("K132").Value = 3622.95
$ws.Range("M132").Value = -1092.95

$ws = $wb.Worksheets.Item("LTW")
# LTW row 68
$ws.Range("H68").Value = 3593.2
$ws.Range("I68").Value = 2992.3333
$ws.Range("K68").Value = 2992.3333
$ws.Range("M68").Value = -2243.3333

# LTW row 71
$ws.Range("H71").Value = 3593.2
$ws.Range("I71").Value = 2992.3333
$ws.Range("K71").Value = 14961.6665
$ws.Range("M71").Value = -11217.6665

# LTW row 82
$ws.Range("H82").Value = 1449.5
$ws.Range("J82").Value = 1449.5
$ws.Range("L82").Value = 1449.5
$ws.Range("N82").Value = -2171.5

# LTW row 85
$ws.Range("H85").Value = 1449.5
$ws.Range("J85").Value = 1449.5
$ws.Range("L85").Value = 1449.5
$ws.Range("N85").Value = -3945.5

# LTW row 93
$ws.Range("H93").Value = 3393.5
$ws.Range("I93").Value = 3528.6667
$ws.Range("J93").Value = 2988
$ws.Range("K93").Value = 3528.6667
$ws.Range("L93").Value = 2988
$ws.Range("M93").Value = -2280.6667
$ws.Range("N93").Value = -5484

# LTW row 100
$ws.Range("H100").Value = 669.3333
$ws.Range("I100").Value = 669.3333
$ws.Range("K100").Value = 669.3333
$ws.Range("M100").Value = -128.3333

# LTW row 122
$ws.Range("H122").Value = 5066.1665
$ws.Range("I122").Value = 4679.4
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 14038.2
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -11588.2
$ws.Range("N122").Value = -25900

$ws = $wb.Worksheets.Item("WVR")
# WVR row 81
$ws.Range("H81").Value = 10744.363
$ws.Range("J81").Value = 13833.833
$ws.Range("L81").Value = 27667.666
$ws.Range("N81").Value = -29789.666

# WVR row 84
$ws.Range("H84").Value = 10744.363
$ws.Range("J84").Value = 13833.833
$ws.Range("L84").Value = 138338.33
$ws.Range("N84").Value = -148946.33
